$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 4

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 5

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 6

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 3

# Update column widths (target OOXML widths ~18.2 and ~13.47 characters;
# the host engine rounds ColumnWidth to a 1/6-character pixel grid, so we
# pick the ColumnWidth input that lands closest to those target widths)
$ws.Columns.Item(1).ColumnWidth = 17.3
$ws.Columns.Item(2).ColumnWidth = 12.65

# Update selection
$ws.Range("A5").Select()
